# The sheet holds one price-quote row per (date, quality) observation for
# "Plátano" at "Vega Monumental Concepción". This edit inserts one missing
# weekly observation before the existing row 403, pushing every following
# row down by one (dimension grows from A1:T467 to A1:T468).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 403..467 down to 404..468, leaving a blank row 403 that
# inherits the formatting (incl. the date-style on column D) from the row
# that used to be there.
$ws.Rows.Item(403).Insert()

# Populate the newly inserted row 403 with the new weekly observation.
$ws.Cells.Item(403, 1).Value = 11
$ws.Cells.Item(403, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(403, 3).Value = "Bíobío"
$ws.Cells.Item(403, 4).Value = 44637
$ws.Cells.Item(403, 5).Value = 8
$ws.Cells.Item(403, 6).Value = "Fruta"
$ws.Cells.Item(403, 7).Value = 100108
$ws.Cells.Item(403, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(403, 9).Value = 100108006
$ws.Cells.Item(403, 10).Value = "Plátano"
$ws.Cells.Item(403, 11).Value = "Sin especificar"
$ws.Cells.Item(403, 12).Value = "Pintón"
$ws.Cells.Item(403, 13).Value = 350
$ws.Cells.Item(403, 14).Value = 19000
$ws.Cells.Item(403, 15).Value = 20000
$ws.Cells.Item(403, 16).Value = 19571
$ws.Cells.Item(403, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(403, 18).Value = "Ecuador"
$ws.Cells.Item(403, 19).Value = 979
$ws.Cells.Item(403, 20).Value = 20
